$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new player row (Pete Sampras) as row 10, mirroring existing data layout
$ws.Range("A10").Value = 13
$ws.Range("B10").Value = "Pete"
$ws.Range("C10").Value = "Sampras"

# Update the active selection to D10 (matches the post-edit selection saved with the file)
$ws.Range("D10").Select()
